# close #37: Check if the indicator is not level zero or negative
#
# - B3 (nivel/"level" indicator): 2 -> 0
# - B5 (nivel/"level" indicator): 3 -> -1
# - Row heights re-flow for rows 4,5,6,7,8,9 (wrapped text grows to fit
#   the now-longer "level" values rendered in the sheet)
# - Font color for the body font (used by almost every data cell) is
#   pinned to explicit black (FF000000) instead of the theme-1 reference

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- data fixes -----------------------------------------------------
$ws.Range("B3").Value = 0
$ws.Range("B5").Value = -1

# --- row heights ------------------------------------------------------
$ws.Rows.Item(4).RowHeight = 73.5
$ws.Rows.Item(5).RowHeight = 100.5
$ws.Rows.Item(6).RowHeight = 73.5
$ws.Rows.Item(7).RowHeight = 73.5
$ws.Rows.Item(8).RowHeight = 100.5
$ws.Rows.Item(9).RowHeight = 60

# --- font color -------------------------------------------------------
# Pin the body font (every non-header cell, A2:Y866, plus the trailing
# blank header cells L1:Y1 that already share that font) to an explicit
# black instead of the "automatic"/theme-1 color.
$ws.Range("A2:Y866").Font.Color = 0
$ws.Range("L1:Y1").Font.Color = 0
